# Apply cell-level updates to rows 2-9 per the source diff (rotation of record data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 96905066
$ws.Range("B2").Value = 77959
$ws.Range("D2").Value = "NT"
$ws.Range("E2").Value = 1797
$ws.Range("F2").Value = "Mjölig dropplav"
$ws.Range("G2").Value = "Cliostomum leprosum"
$ws.Range("H2").Value = "(Räsänen) Holien & Tønsberg"
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()
$ws.Range("P2").Value = "Östansbo, Hundbäcken, Dlr"
$ws.Range("Q2").Value = 509096.042364086
$ws.Range("R2").Value = 6663775.208665404
$ws.Range("Y2").Value = "'2021-10-30"
$ws.Range("AA2").Value = "'2021-10-30"
$ws.Range("AI2").Value = "inslag av tall"
$ws.Range("AL2").Value = "grov gran"
$ws.Range("AO2").Value = "grov gran"
$ws.Range("A3").Value = 16798841
$ws.Range("B3").Value = 88853
$ws.Range("E3").Value = 4189
$ws.Range("F3").Value = "Kamjordstjärna"
$ws.Range("G3").Value = "Geastrum pectinatum"
$ws.Range("H3").Value = "Pers.:Pers."
$ws.Range("I3").Value = "'6"
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("Q3").Value = 509176.1411202319
$ws.Range("R3").Value = 6663496.209875446
$ws.Range("Y3").Value = "'2014-10-30"
$ws.Range("AA3").Value = "'2014-10-30"
$ws.Range("AH3").Value = "Lågörtgranskog"
$ws.Range("AI3").Value = "100-årig skog"
$ws.Range("AJ3").ClearContents()
$ws.Range("AK3").ClearContents()
$ws.Range("AO3").Value = "gammal myrstack"
$ws.Range("A4").Value = 74249883
$ws.Range("B4").Value = 89356
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5447
$ws.Range("F4").Value = "Vedticka"
$ws.Range("G4").Value = "Fuscoporia viticola"
$ws.Range("H4").Value = "(Schwein.) Murrill"
$ws.Range("Q4").Value = 509213.9838288009
$ws.Range("R4").Value = 6663342.23806499
$ws.Range("AO4").Value = "gammal granlåga # Picea abies"
$ws.Range("A5").Value = 74249894
$ws.Range("B5").Value = 89851
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5467
$ws.Range("F5").Value = "Kådvaxskinn"
$ws.Range("G5").Value = "Phlebia serialis"
$ws.Range("H5").Value = "(Fr.:Fr.) Donk"
$ws.Range("P5").Value = "Östansbo vid Hundbäcken, Dlr"
$ws.Range("Q5").Value = 509150.8558177771
$ws.Range("R5").Value = 6663421.85664573
$ws.Range("Y5").Value = "'2018-11-14"
$ws.Range("AA5").Value = "'2018-11-14"
$ws.Range("AH5").Value = "Blåbärsgranskog"
$ws.Range("AI5").Value = "mossigt"
$ws.Range("AO5").Value = "gammal klen granlåga # Picea abies"
$ws.Range("A6").Value = 96640115
$ws.Range("B6").Value = 90319
$ws.Range("E6").Value = 4769
$ws.Range("F6").Value = "Svavelriska"
$ws.Range("G6").Value = "Lactarius scrobiculatus"
$ws.Range("H6").Value = "(Scop.:Fr.) Fr."
$ws.Range("A7").Value = 96640119
$ws.Range("B7").Value = 90674
$ws.Range("D7").Value = "LC"
$ws.Range("E7").Value = 5964
$ws.Range("F7").Value = "Fjällig taggsvamp s.str."
$ws.Range("G7").Value = "Sarcodon imbricatus s.str."
$ws.Range("H7").Value = "(L.:Fr.) P.Karst."
$ws.Range("P7").Value = "Hundbäcken, Östansbo, Dlr"
$ws.Range("Q7").Value = 509116.3226223197
$ws.Range("R7").Value = 6663249.75091619
$ws.Range("Y7").Value = "'2021-10-13"
$ws.Range("AA7").Value = "'2021-10-13"
$ws.Range("AJ7").Value = "gran"
$ws.Range("AK7").Value = "Picea abies"
$ws.Range("AO7").Value = "Picea abies"
$ws.Range("A8").Value = 96905129
$ws.Range("B8").Value = 77959
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 1797
$ws.Range("F8").Value = "Mjölig dropplav"
$ws.Range("G8").Value = "Cliostomum leprosum"
$ws.Range("H8").Value = "(Räsänen) Holien & Tønsberg"
$ws.Range("Q8").Value = 509157.4464191007
$ws.Range("R8").Value = 6663582.916298064
$ws.Range("AO8").ClearContents()
$ws.Range("A9").Value = 96905141
$ws.Range("B9").Value = 89832
$ws.Range("D9").Value = "VU"
$ws.Range("E9").Value = 1209
$ws.Range("F9").Value = "Rynkskinn"
$ws.Range("G9").Value = "Phlebia centrifuga"
$ws.Range("H9").Value = "P.Karst."
$ws.Range("Q9").Value = 509087.8970259681
$ws.Range("R9").Value = 6663434.16527954
$ws.Range("AI9").ClearContents()
$ws.Range("AL9").ClearContents()
$ws.Range("AO9").Value = "barklös granlåga"
